$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "99.008.53"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "3.288.81"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'254.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "'625.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("D7").Value = "'1.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +23.84%  "
$ws.Range("D8").Value = "'0.401"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.67%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'0.979"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +23.58%  "
$ws.Range("D11").Value = "3.286.19"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "'0.203"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").Value = "'40.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.26%  "
$ws.Range("D14").Value = "98.760.25"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "3.897.87"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "'5.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "3.291.79"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("E19").Value = "  -3.68%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'15.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.93%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "'6.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.32%  "
$ws.Range("D22").Value = "'489.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.22%  "
$ws.Range("D23").Value = "'9.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").Value = "'0.344"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +39.69%  "
$ws.Range("D26").Value = "'5.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").Value = "'89.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").Value = "'12.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").Value = "3.461.37"
$ws.Range("E29").Value = "  -3.07%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.145"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.51%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'0.189"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").Value = "'10.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.84%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'27.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("E36").Value = "  +7.35%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'7.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'1.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "'492.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.77%  "
$ws.Range("D41").Value = "'24.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'3.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.28%  "
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'0.777"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("D47").Value = "'158.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'4.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.84%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'7.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.62%  "
$ws.Range("D51").Value = "'0.850"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.49%  "
